# Configure project for selenium support (local chrome driver version 87)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("project setup"): Estimated/Actual columns updated
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1

# Row 4 ("cofigure selenium and chrome driver for local-development"):
# fill in Actual Time and % complete
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D4").NumberFormat = "0%"

# Update the active selection to reflect where the edit left off
$ws.Range("E4").Select()
